$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F ("想去人数")
$updates = @{
    3  = 54
    5  = 36
    7  = 2756
    9  = 1770
    11 = 79
    12 = 617
    15 = 152
}

# Both "展览" and "全部类型" sheets contain the same event data and need the same update
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
